$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.460.51"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "3.620.36"
$ws.Range("E3").Value = "  +2.14%  "
$ws.Range("D5").Value = "'605.52"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "'196.27"
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("D7").Value = "'0.627"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.206"
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").Value = "'53.66"
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").Value = "4.189.18"
$ws.Range("E14").Value = "  +1.86%  "
$ws.Range("D15").Value = "'13.01"
$ws.Range("E15").Value = "  +2.33%  "
$ws.Range("D16").Value = "'595.47"
$ws.Range("E16").Value = "  -1.52%  "
$ws.Range("D17").Value = "70.519.59"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").Value = "3.614.93"
$ws.Range("E18").Value = "  +2.21%  "
$ws.Range("E19").Value = "  -1.16%  "
$ws.Range("D21").Value = "'0.997"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "'17.82"
$ws.Range("E22").Value = "  -1.67%  "
$ws.Range("D23").Value = "'5.18"
$ws.Range("E23").Value = "  -1.94%  "
$ws.Range("D24").Value = "'101.82"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("E26").Value = "  -3.67%  "
$ws.Range("E27").Value = "  -1.85%  "
$ws.Range("D28").Value = "'9.62"
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D29").Value = "'33.83"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "'4.69"
$ws.Range("E30").Value = "  +7.11%  "
$ws.Range("D31").Value = "'7.23"
$ws.Range("E31").Value = "  +1.15%  "
$ws.Range("D32").Value = "'12.31"
$ws.Range("E32").Value = "  -2.86%  "
$ws.Range("D33").Value = "'0.118"
$ws.Range("E33").Value = "  +1.71%  "
$ws.Range("D34").Value = "'63.67"
$ws.Range("E34").Value = "  +0.43%  "
$ws.Range("D35").Value = "0.0₃0896"
$ws.Range("E35").Value = "  +6.01%  "
$ws.Range("D36").Value = "3.907.36"
$ws.Range("E36").Value = "  +3.14%  "
$ws.Range("D37").Value = "'541.27"
$ws.Range("E37").Value = "  +10.56%  "
$ws.Range("D38").Value = "'3.14"
$ws.Range("E38").Value = "  +1.94%  "
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("D40").Value = "'37.02"
$ws.Range("E40").Value = "  +0.51%  "
$ws.Range("E41").Value = "  -1.46%  "
$ws.Range("D42").Value = "'3.52"
$ws.Range("E42").Value = "  -4.77%  "
$ws.Range("E43").Value = "  -0.83%  "
$ws.Range("D44").Value = "'0.0456"
$ws.Range("E44").Value = "  -0.54%  "
$ws.Range("D45").Value = "'3.42"
$ws.Range("E45").Value = "  +3.77%  "
$ws.Range("D46").Value = "'2.86"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("E48").Value = "  -0.97%  "
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").Value = "'0.000253"
$ws.Range("E50").Value = "  +0.83%  "
$ws.Range("E51").Value = "  +0.45%  "
